$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 4.45 = 17310.19 pesos`n✅ 17310.19 pesos = 4.41 = 941.73 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the rate cells N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 224.8
$wsTasas.Range("O10").Value = 3891.33
$wsTasas.Range("N12").Value = 3929
$wsTasas.Range("O12").Value = 213.75
